$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell -> new text value (prices and percentages stored as literal text)
$updates = [ordered]@{
    "D2" = "279.63"
    "E2" = "5.80%"
    "D3" = "27.06"
    "E3" = "1.49%"
    "D4" = "4.939"
    "E4" = "5.10%"
    "D5" = "0.06374"
    "E5" = "4.26%"
    "D6" = "6.954"
    "E6" = "3.24%"
    "D7" = "3.356"
    "E7" = "6.08%"
    "D8" = "0.8839"
    "E8" = "3.86%"
    "D9" = "0.9454"
    "E9" = "3.89%"
    "D10" = "0.1469"
    "E10" = "4.26%"
    "D11" = "0.05138"
    "E11" = "7.51%"
    "D12" = "0.07397"
    "E12" = "4.24%"
    "D13" = "0.03136"
    "E13" = "0.24%"
    "D14" = "0.09052"
    "E14" = "0.08%"
    "D15" = "0.001557"
    "E15" = "1.29%"
    "D16" = "0.0006288"
    "E16" = "1.94%"
    "D17" = "0.005966"
    "E17" = "-0.12%"
    "D18" = "3.489"
    "E18" = "1.04%"
    "E19" = "7.04%"
    "E20" = "0.85%"
    "D21" = "0.1329"
    "E21" = "3.81%"
    "D22" = "3.889"
    "E22" = "-5.58%"
    "D23" = "0.04326"
    "E23" = "1.97%"
    "D24" = "0.001174"
    "E24" = "-0.39%"
    "D25" = "0.003644"
    "E25" = "-10.41%"
    "D26" = "0.0001198"
    "E26" = "-0.14%"
    "D27" = "0.0001693"
    "E27" = "-12.64%"
    "D40" = "0.04074"
    "E40" = "3.51%"
    "D41" = "0.006621"
    "E41" = "58.62%"
    "D42" = "0.1168"
    "E42" = "4.69%"
    "D43" = "0.002346"
    "E43" = "11.26%"
    "E44" = "7.99%"
    "D45" = "0.00005250"
    "E45" = "3.65%"
    "E46" = "0.05%"
    "D47" = "2.373"
    "E47" = "818.50%"
    "D48" = "0.02258"
    "E48" = "6.45%"
    "E49" = "0.05%"
    "E50" = "-0.02%"
}

foreach ($cell in $updates.Keys) {
    $range = $ws.Range($cell)
    $range.NumberFormat = "@"
    $range.Value = $updates[$cell]
}
